# Generate Report for Archive
#
# The localization status moved from "Ready for handoff" to "In Translation"
# for the a3da39c9-... source file, on every sheet that surfaces that value
# (Overview!E2:F2, zh-cn!C2, de-de!C2 all share the same cell text). After the
# shorter label is written, the Status-ish columns that displayed it are
# narrowed to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1:F2").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1:C2").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1:C2").ColumnWidth = 12.5
